$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r AnalysisUnit_Variable")

$newVars = @(
  "PREAM_UTIL_ACC_L1M",
  "PREAM_MAX_MM_L1M",
  "EDF_ImpliedRating",
  "Bond_ImpliedRating",
  "CDS_ImpliedRating",
  "5Y_CDS_Spread",
  "Distance_to_default",
  "Rating_stress",
  "Notizie_Pregiudizievoli",
  "Colore_Modulo_Dati_di_Mercato",
  "NPAF_IS_FAC_L1M",
  "NPAF_IS_ANT_L1M",
  "NPAF_IS_MAX_FAC_L3M",
  "NPAF_IS_MAX_ANT_L3M",
  "NPAF_IS_M_FAC_L3M",
  "NPAF_IS_M_ANT_L3M",
  "NPAF_IS_MAX_DOC_L3M",
  "NPAF_IS_M_DOC_L3M",
  "NPAF_IS_FAC_FAT_L1M",
  "NPAF_IS_ANT_FAT_L1M",
  "NPAF_IS_MAX_FAC_FAT_L3M",
  "NPAF_IS_MAX_ANT_FAT_L3M",
  "NPAF_IS_M_FAC_FAT_L3M",
  "NPAF_IS_M_ANT_FAT_L3M",
  "NPAF_IS_FAC_DB_L1M",
  "NPAF_IS_ANT_DB_L1M",
  "NPAF_IS_MAX_FAC_DB_L3M",
  "NPAF_IS_MAX_ANT_DB_L3M",
  "NPAF_IS_M_FAC_DB_L3M",
  "NPAF_IS_M_ANT_DB_L3M",
  "NPAF_IS_FAC_MOL_L1M",
  "NPAF_IS_ANT_MOL_L1M",
  "NPAF_IS_MAX_FAC_MOL_L3M",
  "NPAF_IS_MAX_ANT_MOL_L3M",
  "NPAF_IS_M_FAC_MOL_L3M",
  "NPAF_IS_M_ANT_MOL_L3M",
  "NPAF_IS_MAX_DOC_DB_L3M",
  "NPAF_IS_M_DOC_DB_L3M",
  "NPAF_IS_ANT_ACC_L1M",
  "ANTEXP_SCAD_ACC_1G_L1M",
  "ANTEXP_SCAD_ACC_60G_L1M",
  "ANTEXP_SCAD_ACC_1G_M_L3M",
  "ANTEXP_SCAD_ACC_30G_M_L3M",
  "ANTEXP_SCAD_ACC_60G_M_L3M",
  "FINIMP_SCAD_ACC_1G_L1M",
  "FINIMP_SCAD_ACC_60G_L1M",
  "FINIMP_SCAD_ACC_1G_M_L3M",
  "FINIMP_SCAD_ACC_30G_M_L3M",
  "FINIMP_SCAD_ACC_60G_M_L3M"
)

$startRow = 100
for ($i = 0; $i -lt $newVars.Count; $i++) {
  $row = $startRow + $i
  $name = $newVars[$i]
  $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
  $ws.Cells.Item($row, 2).Value = $name
  $ws.Cells.Item($row, 3).Value = $name
  $ws.Cells.Item($row, 5).Value = "CUSTOMER"
  $ws.Cells.Item($row, 6).Value = $name
}

$ws.Activate()
$ws.Range("B147").Select()
